$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns remain plain text (they contain values like
# "1.007" or "26.168.73" that must not be reinterpreted as numbers/dates).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.981.56"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "1.800.54"
$ws.Range("E3").Value = "  -2.07%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.72%  "
$ws.Range("D5").Value = "239.30"
$ws.Range("E5").Value = "  -7.92%  "
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").Value = "0.5058"
$ws.Range("E7").Value = "  -3.51%  "
$ws.Range("D8").Value = "0.2464"
$ws.Range("E8").Value = "  -22.65%  "
$ws.Range("D9").Value = "0.06067"
$ws.Range("E9").Value = "  -10.56%  "
$ws.Range("D10").Value = "1.821.99"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "0.06860"
$ws.Range("E11").Value = "  -11.44%  "
$ws.Range("D12").Value = "14.87"
$ws.Range("E12").Value = "  -20.57%  "
$ws.Range("D13").Value = "78.96"
$ws.Range("E13").Value = "  -10.04%  "
$ws.Range("D14").Value = "0.5975"
$ws.Range("E14").Value = "  -23.73%  "
$ws.Range("D15").Value = "4.406"
$ws.Range("E15").Value = "  -12.02%  "
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "26.049.88"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").Value = "11.28"
$ws.Range("E19").Value = "  -18.42%  "
$ws.Range("D20").Value = "2.065.25"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "0.000005856"
$ws.Range("E21").Value = "  -26.30%  "
$ws.Range("D22").Value = "3.962"
$ws.Range("E22").Value = "  -14.26%  "
$ws.Range("D23").Value = "5.167"
$ws.Range("E23").Value = "  -13.31%  "
$ws.Range("D24").Value = "7.929"
$ws.Range("E24").Value = "  -15.21%  "
$ws.Range("D25").Value = "130.63"
$ws.Range("E25").Value = "  -7.81%  "
$ws.Range("D26").Value = "1.864"
$ws.Range("E26").Value = "  -14.42%  "
$ws.Range("D27").Value = "14.51"
$ws.Range("E27").Value = "  -14.07%  "
$ws.Range("D28").Value = "98.89"
$ws.Range("E28").Value = "  -11.21%  "
$ws.Range("D29").Value = "1.215"
$ws.Range("E29").Value = "  -27.61%  "
$ws.Range("D30").Value = "0.08205"
$ws.Range("E30").Value = "  -5.60%  "
$ws.Range("D31").Value = "3.645"
$ws.Range("E31").Value = "  -12.29%  "
$ws.Range("D32").Value = "2.759"
$ws.Range("E32").Value = "  -4.03%  "
$ws.Range("D33").Value = "3.175"
$ws.Range("E33").Value = "  -21.94%  "
$ws.Range("D34").Value = "0.04262"
$ws.Range("E34").Value = "  -12.62%  "
$ws.Range("D35").Value = "1.042"
$ws.Range("E35").Value = "  -8.07%  "
$ws.Range("D36").Value = "2.901"
$ws.Range("E36").Value = "  -6.21%  "
$ws.Range("D37").Value = "0.6236"
$ws.Range("E37").Value = "  -14.24%  "
$ws.Range("D38").Value = "2.072"
$ws.Range("E38").Value = "  -7.29%  "
$ws.Range("D39").Value = "1.007"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").Value = "0.8146"
$ws.Range("E40").Value = "  -8.68%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01454"
$ws.Range("E41").Value = "  -17.03%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "98.96"
$ws.Range("E42").Value = "  -9.64%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.3831"
$ws.Range("E43").Value = "  -19.61%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.151"
$ws.Range("E44").Value = "  -13.20%  "
$ws.Range("D45").Value = "0.05264"
$ws.Range("E45").Value = "  -9.99%  "
$ws.Range("D46").Value = "6.190"
$ws.Range("E46").Value = "  -19.13%  "
$ws.Range("D47").Value = "53.14"
$ws.Range("E47").Value = "  -10.65%  "
$ws.Range("B48").Value = "USDD"
$ws.Range("C48").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D48").Value = "1.009"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.1032"
$ws.Range("E49").Value = "  -16.10%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "29.21"
$ws.Range("E50").Value = "  -16.11%  "
$ws.Range("B51").Value = "TrueUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  +0.34%  "

# Restore the default cell style (NumberFormat=@ above creates a temporary style)
$ws.Range("B2:E51").Style = "Normal"
